# Append one new data row (row 50) to each of the four sensor-log sheets,
# duplicating the previous last row (row 49) but with the timestamp hour
# advanced by one (08:xx:xx -> 09:xx:xx) and the ID_DEC value updated.
#
# Column layout: A=time, B=总长, C=ID, D=实际长度, E=和校验,
#                F=总长_DEC, G=ID_DEC, H=实际长度_DEC, I=和校验_DEC

$wb = $excel.ActiveWorkbook

$newRows = @{
  "ROW35-FE-LIFTER"  = @{
    A = "2025-03-06 09:42:06"
    B = "0x01,0x90 "
    C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"
    D = "0x01,0x90,"
    E = "0x d"
    F = 400
    G = "568631262647113770877196"
    H = 400
    I = 13
  }
  "ROW35-MID-LIFTER" = @{
    A = "2025-03-06 09:29:35"
    B = "0x01,0x90 "
    C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"
    D = "0x01,0x90,"
    E = "0x e"
    F = 400
    G = "568631262647113770942732"
    H = 400
    I = 14
  }
  "ROW02-FE-LIFTER"  = @{
    A = "2025-03-06 09:51:45"
    B = "0x01,0x90 "
    C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"
    D = "0x01,0x90,"
    E = "0xff"
    F = 400
    G = "568631262647113769959692"
    H = 400
    I = 255
  }
  "ROW02-MID-LIFTER" = @{
    A = "2025-03-06 09:41:15"
    B = "0x01,0x90 "
    C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
    D = "0x01,0x90,"
    E = "0x 3"
    F = 400
    G = "568631262647113769959692"
    H = 400
    I = 3
  }
}

foreach ($sheetName in $newRows.Keys) {
  $ws = $wb.Worksheets.Item($sheetName)
  $row = $newRows[$sheetName]
  $targetRow = 50

  # Text columns: force text storage (avoids numeric/date auto-conversion,
  # which matters most for the long ID_DEC digit string in column G), then
  # drop back to the default "Normal" style so the cell ends up unstyled
  # (matching the rest of the sheet's plain data cells) while the stored
  # value stays textual.
  foreach ($col in @(1, 2, 3, 4, 5, 7)) {
    $letter = [char](64 + $col)
    $cell = $ws.Cells.Item($targetRow, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $row[$letter]
    $cell.Style = "Normal"
  }

  # Numeric columns.
  $ws.Cells.Item($targetRow, 6).Value = $row["F"]
  $ws.Cells.Item($targetRow, 8).Value = $row["H"]
  $ws.Cells.Item($targetRow, 9).Value = $row["I"]
}
